$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1835.72
$ws.Range("I43").Value = 426.66666
$ws.Range("J43").Value = 2628.3125
$ws.Range("K43").Value = 426.66666
$ws.Range("L43").Value = 2628.3125
$ws.Range("M43").Value = -357.66666
$ws.Range("N43").Value = -2766.3125
$ws.Range("H47").Value = 28383.334
$ws.Range("I47").Value = 20575
$ws.Range("J47").Value = 44000
$ws.Range("K47").Value = 20575
$ws.Range("L47").Value = 44000
$ws.Range("M47").Value = -19603
$ws.Range("N47").Value = -45944
$ws.Range("H137").Value = 748574.4
$ws.Range("I137").Value = 1245082.1
$ws.Range("J137").Value = 3812.875
$ws.Range("K137").Value = 3735246.3
$ws.Range("L137").Value = 11438.625
$ws.Range("M137").Value = -3732696.3
$ws.Range("N137").Value = -16538.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1415.8334
$ws.Range("I26").Value = 1415.8334
$ws.Range("K26").Value = 1415.8334
$ws.Range("M26").Value = -1085.8334
$ws.Range("H61").Value = 2689
$ws.Range("I61").Value = 1501.875
$ws.Range("K61").Value = 1501.875
$ws.Range("M61").Value = -1289.875
$ws.Range("H74").Value = 3132.5
$ws.Range("I74").Value = 702.5
$ws.Range("J74").Value = 5388.9287
$ws.Range("K74").Value = 702.5
$ws.Range("L74").Value = 5388.9287
$ws.Range("M74").Value = 171.5
$ws.Range("N74").Value = -7136.9287
$ws.Range("H77").Value = 3132.5
$ws.Range("I77").Value = 702.5
$ws.Range("J77").Value = 5388.9287
$ws.Range("K77").Value = 3512.5
$ws.Range("L77").Value = 26944.6435
$ws.Range("M77").Value = 855.5
$ws.Range("N77").Value = -35680.64350000001
$ws.Range("H136").Value = 2689
$ws.Range("I136").Value = 1501.875
$ws.Range("K136").Value = 4505.625
$ws.Range("M136").Value = -1955.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1960.4828
$ws.Range("I86").Value = 1681.4166
$ws.Range("J86").Value = 3300
$ws.Range("K86").Value = 1681.4166
$ws.Range("L86").Value = 3300
$ws.Range("M86").Value = -558.4166
$ws.Range("N86").Value = -5546
$ws.Range("H89").Value = 1960.4828
$ws.Range("I89").Value = 1681.4166
$ws.Range("J89").Value = 3300
$ws.Range("K89").Value = 8407.083000000001
$ws.Range("L89").Value = 16500
$ws.Range("M89").Value = -2791.083000000001
$ws.Range("N89").Value = -27732
$ws.Range("H94").Value = 851.64
$ws.Range("I94").Value = 845.13635
$ws.Range("J94").Value = 899.3333
$ws.Range("K94").Value = 845.13635
$ws.Range("L94").Value = 899.3333
$ws.Range("M94").Value = -394.13635
$ws.Range("N94").Value = -1801.3333
$ws.Range("H107").Value = 372141.44
$ws.Range("I107").Value = 487064.22
$ws.Range("J107").Value = 1834.7778
$ws.Range("K107").Value = 487064.22
$ws.Range("L107").Value = 1834.7778
$ws.Range("M107").Value = -485144.22
$ws.Range("N107").Value = -5674.7778
$ws.Range("H134").Value = 33783.89
$ws.Range("I134").Value = 42730.207
$ws.Range("K134").Value = 128190.621
$ws.Range("M134").Value = -125655.621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1658.5
$ws.Range("I31").Value = 997.6923
$ws.Range("J31").Value = 2885.7144
$ws.Range("K31").Value = 997.6923
$ws.Range("L31").Value = 2885.7144
$ws.Range("M31").Value = -702.6923
$ws.Range("N31").Value = -3475.7144
$ws.Range("H34").Value = 1658.5
$ws.Range("I34").Value = 997.6923
$ws.Range("J34").Value = 2885.7144
$ws.Range("K34").Value = 997.6923
$ws.Range("L34").Value = 2885.7144
$ws.Range("M34").Value = -795.6923
$ws.Range("N34").Value = -3289.7144
$ws.Range("H58").Value = 3093.513
$ws.Range("I58").Value = 1370.45
$ws.Range("J58").Value = 4907.263
$ws.Range("K58").Value = 1370.45
$ws.Range("L58").Value = 4907.263
$ws.Range("M58").Value = -1167.45
$ws.Range("N58").Value = -5313.263
$ws.Range("H98").Value = 38990
$ws.Range("J98").Value = 38990
$ws.Range("L98").Value = 38990
$ws.Range("N98").Value = -43482
$ws.Range("H107").Value = 3788357
$ws.Range("I107").Value = 5952826
$ws.Range("J107").Value = 536.5
$ws.Range("K107").Value = 5952826
$ws.Range("L107").Value = 536.5
$ws.Range("M107").Value = -5950906
$ws.Range("N107").Value = -4376.5
$ws.Range("H132").Value = 2166.8438
$ws.Range("I132").Value = 1927.5217
$ws.Range("K132").Value = 5782.5651
$ws.Range("M132").Value = -3252.5651
$ws.Range("H134").Value = 2850.2593
$ws.Range("I134").Value = 3044.913
$ws.Range("J134").Value = 1731
$ws.Range("K134").Value = 9134.739
$ws.Range("L134").Value = 5193
$ws.Range("M134").Value = -6599.739
$ws.Range("N134").Value = -10263
$ws.Range("H136").Value = 3093.513
$ws.Range("I136").Value = 1370.45
$ws.Range("J136").Value = 4907.263
$ws.Range("K136").Value = 4111.35
$ws.Range("L136").Value = 14721.789
$ws.Range("M136").Value = -1561.35
$ws.Range("N136").Value = -19821.789
$ws.Range("H137").Value = 52999.75
$ws.Range("J137").Value = 52999.75
$ws.Range("L137").Value = 52999.75
$ws.Range("N137").Value = -63199.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2666.6667
$ws.Range("J35").Value = 2666.6667
$ws.Range("L35").Value = 8000.000100000001
$ws.Range("N35").Value = -8576.000100000001
$ws.Range("H57").Value = 1642.8572
$ws.Range("J57").Value = 1750
$ws.Range("L57").Value = 5250
$ws.Range("N57").Value = -6368
$ws.Range("H74").Value = 2000
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -11122
$ws.Range("H77").Value = 2000
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 27000
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -37608
$ws.Range("H88").Value = 834277.5600000001
$ws.Range("J88").Value = 999.7646999999999
$ws.Range("L88").Value = 2999.2941
$ws.Range("N88").Value = -3855.2941
$ws.Range("H91").Value = 834277.5600000001
$ws.Range("J91").Value = 999.7646999999999
$ws.Range("L91").Value = 2999.2941
$ws.Range("N91").Value = -5963.2941
$ws.Range("H99").Value = 2762.3333
$ws.Range("I99").Value = 1589
$ws.Range("J99").Value = 3349
$ws.Range("K99").Value = 4767
$ws.Range("L99").Value = 10047
$ws.Range("M99").Value = -2521
$ws.Range("N99").Value = -14539
$ws.Range("H102").Value = 5414.2856
$ws.Range("J102").Value = 7180
$ws.Range("L102").Value = 21540
$ws.Range("N102").Value = -26408
$ws.Range("H113").Value = 523.86664
$ws.Range("I113").Value = 511.7857
$ws.Range("J113").Value = 529.3226
$ws.Range("K113").Value = 1535.3571
$ws.Range("L113").Value = 1587.9678
$ws.Range("M113").Value = 634.6428999999998
$ws.Range("N113").Value = -5927.9678
$ws.Range("H131").Value = 1615484.2
$ws.Range("I131").Value = 7780
$ws.Range("J131").Value = 1853662.6
$ws.Range("K131").Value = 23340
$ws.Range("L131").Value = 5560987.800000001
$ws.Range("M131").Value = -18300
$ws.Range("N131").Value = -5571067.800000001
$ws.Range("H133").Value = 6088.1514
$ws.Range("I133").Value = 3032
$ws.Range("J133").Value = 7416.913
$ws.Range("K133").Value = 9096
$ws.Range("L133").Value = 22250.739
$ws.Range("M133").Value = -4036
$ws.Range("N133").Value = -32370.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 7074.364
$ws.Range("I36").Value = 10017
$ws.Range("J36").Value = 5392.857
$ws.Range("K36").Value = 10017
$ws.Range("L36").Value = 5392.857
$ws.Range("M36").Value = -9532
$ws.Range("N36").Value = -6362.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 14649
$ws.Range("I45").Value = 14500
$ws.Range("J45").Value = 14698.667
$ws.Range("K45").Value = 14500
$ws.Range("L45").Value = 14698.667
$ws.Range("M45").Value = -14093
$ws.Range("N45").Value = -15512.667
$ws.Range("H62").Value = 19899.4
$ws.Range("I62").Value = 15000
$ws.Range("J62").Value = 21124.25
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 21124.25
$ws.Range("M62").Value = -14376
$ws.Range("N62").Value = -22372.25
$ws.Range("H65").Value = 19899.4
$ws.Range("I65").Value = 15000
$ws.Range("J65").Value = 21124.25
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 63372.75
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -69612.75
$ws.Range("H76").Value = 6490
$ws.Range("I76").Value = 4980
$ws.Range("K76").Value = 4980
$ws.Range("M76").Value = -4642
$ws.Range("H79").Value = 6490
$ws.Range("I79").Value = 4980
$ws.Range("K79").Value = 4980
$ws.Range("M79").Value = -3810
$ws.Range("H100").Value = 1434.6666
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 1152
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 1152
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -2234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 281.8095
$ws.Range("I100").Value = 258.8421
$ws.Range("K100").Value = 517.6842
$ws.Range("M100").Value = 23.31579999999997
$ws.Range("H107").Value = 425.67743
$ws.Range("I107").Value = 268
$ws.Range("K107").Value = 804
$ws.Range("M107").Value = 1116
$ws.Range("H126").Value = 2200.7144
$ws.Range("I126").Value = 2254.6155
$ws.Range("K126").Value = 6763.8465
$ws.Range("M126").Value = -4293.8465
$ws.Range("H132").Value = 4659.1934
$ws.Range("I132").Value = 6705.4116
$ws.Range("J132").Value = 2174.5
$ws.Range("K132").Value = 20116.2348
$ws.Range("L132").Value = 6523.5
$ws.Range("M132").Value = -17586.2348
$ws.Range("N132").Value = -11583.5
$ws.Range("H136").Value = 8828.120000000001
$ws.Range("I136").Value = 12478.765
$ws.Range("J136").Value = 1070.5
$ws.Range("K136").Value = 37436.295
$ws.Range("L136").Value = 3211.5
$ws.Range("M136").Value = -34886.295
$ws.Range("N136").Value = -8311.5
